$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data values between row 2 and row 4 for columns D, M, N, O, P, S
# (this reflects the weekly data being re-ordered/swapped for the two dates)

$row2D = $ws.Range("D2").Value2
$row2M = $ws.Range("M2").Value2
$row2N = $ws.Range("N2").Value2
$row2O = $ws.Range("O2").Value2
$row2P = $ws.Range("P2").Value2
$row2S = $ws.Range("S2").Value2

$row4D = $ws.Range("D4").Value2
$row4M = $ws.Range("M4").Value2
$row4N = $ws.Range("N4").Value2
$row4O = $ws.Range("O4").Value2
$row4P = $ws.Range("P4").Value2
$row4S = $ws.Range("S4").Value2

$ws.Range("D2").Value2 = $row4D
$ws.Range("M2").Value2 = $row4M
$ws.Range("N2").Value2 = $row4N
$ws.Range("O2").Value2 = $row4O
$ws.Range("P2").Value2 = $row4P
$ws.Range("S2").Value2 = $row4S

$ws.Range("D4").Value2 = $row2D
$ws.Range("M4").Value2 = $row2M
$ws.Range("N4").Value2 = $row2N
$ws.Range("O4").Value2 = $row2O
$ws.Range("P4").Value2 = $row2P
$ws.Range("S4").Value2 = $row2S
